$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("Q2").Value = 3.1
$ws.Range("R2").Value = 1.36
$ws.Range("M3").Value = 1.1
$ws.Range("N3").Value = 7
$ws.Range("Q3").Value = 2.5
$ws.Range("R3").Value = 1.5
$ws.Range("J5").Value = 6.5
$ws.Range("K5").Value = 2.1
$ws.Range("U5").Value = 2.2
$ws.Range("V5").Value = 1.62
$ws.Range("Y5").Value = 21
$ws.Range("AG5").Value = 5.5
$ws.Range("AI5").Value = 9
$ws.Range("AN5").Value = 7.5
$ws.Range("O6").Value = 1.44
$ws.Range("P6").Value = 2.75
$ws.Range("N8").Value = 10
$ws.Range("G11").Value = 2.35
$ws.Range("H11").Value = 3.75
$ws.Range("I11").Value = 2.88
$ws.Range("J11").Value = 2.88
$ws.Range("L11").Value = 3.4
$ws.Range("U11").Value = 1.57
$ws.Range("V11").Value = 2.25
$ws.Range("AL11").Value = 26
$ws.Range("AP11").Value = 19
$ws.Range("AS11").Value = 101
$ws.Range("G13").Value = 2.47
$ws.Range("I13").Value = 2.72
$ws.Range("J13").Value = 3.15
$ws.Range("L13").Value = 3.35
$ws.Range("N13").Value = 6.3
$ws.Range("O13").Value = 1.4
$ws.Range("P13").Value = 2.75
$ws.Range("Q13").Value = 2.2
$ws.Range("R13").Value = 1.62
$ws.Range("V13").Value = 1.82
$ws.Range("W13").Value = 7.1
$ws.Range("X13").Value = 11.5
$ws.Range("Y13").Value = 9.75
$ws.Range("Z13").Value = 27
$ws.Range("AA13").Value = 23
$ws.Range("AB13").Value = 37
$ws.Range("AC13").Value = 6.3
$ws.Range("AG13").Value = 7.7
$ws.Range("AH13").Value = 13
$ws.Range("AI13").Value = 10.25
$ws.Range("AJ13").Value = 32
$ws.Range("AK13").Value = 25
$ws.Range("AL13").Value = 37
$ws.Range("AN13").Value = 4.35
$ws.Range("AO13").Value = 13.5
$ws.Range("AP13").Value = 23
$ws.Range("AQ13").Value = 60
$ws.Range("AR13").Value = 110
$ws.Range("AS13").Value = 350
$ws.Range("AU13").Value = 7.3
$ws.Range("AW13").Value = 4.6
$ws.Range("AX13").Value = 15.5
$ws.Range("AY13").Value = 24
$ws.Range("AZ13").Value = 70
$ws.Range("BA13").Value = 120
$ws.Range("BB13").Value = 350
